$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Lateral (Al) ---
$ws.Range("A2:B2").Font.Color = -16777216   # black / automatic (drop red highlight)
$ws.Range("B2").Value = 100

# --- Row 3: Cobertura (A definir) --- (no change)

# --- Row 4: Chassi (PETG) -> Chassi (PETG) 50% treliça ---
$ws.Range("A4").Value = "Chassi (PETG) 50% treliça"
$ws.Range("B4").Value = 155

# --- Row 5: Suporte motor mov(Al) ---
$ws.Range("A5:B5").Font.Color = -16777216   # black / automatic (drop red highlight)
$ws.Range("B5").Value = 35

# --- Row 6: Mancal arma (Pol) -- add detail columns E/F/G ---
$ws.Range("E6:G6").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E6").Value = "Eixo"
$ws.Range("F6").Value = "Bits"
$ws.Range("G6").Value = "Discos"

# --- Row 7: Arma (Aço + Al) -- turn B7 into SUM of new detail cells ---
$ws.Range("D7").Value = "Detalhes"
$ws.Range("E7:G7").HorizontalAlignment = -4108   # xlCenter
$ws.Range("E7").Value = 47
$ws.Range("F7").Value = 102
$ws.Range("G7").Value = 168
$ws.Range("B7").Formula = "=SUM(E7:G7)"

# --- Row 8: Rolamentos ---
$ws.Range("A8:B8").Font.Color = -16777216   # black / automatic (drop red highlight)
$ws.Range("B8").Value = 10
$ws.Range("D8").Value = "5 cada"

# --- Row 9: Correia + Polias ---
$ws.Range("A9:B9").Font.Color = -16777216   # black / automatic (drop red highlight)
$ws.Range("B9").Value = 55

# --- Row 10: motores mov --- (no change)
# --- Row 11: Rodas --- (no change)

# --- Row 12: motor arma ---
$ws.Range("B12").Value = 57

# --- Row 13: bateria -> bateria (2S 1700mAh); gains red highlight ---
$ws.Range("A13").Value = "bateria (2S 1700mAh)"
$ws.Range("A13:B13").Font.Color = 255   # red
$ws.Range("B13").Value = 80

# --- Row 14: receptor --- (no change)
# --- Row 15: ESC --- (no change)

# --- Row 16: Placa de controle ---
$ws.Range("B16").Value = 40

# --- Row 17: Fios + conectores --- (no change)

# --- Row 18: Parafusos -> Parafusos (Reduzidos) ---
$ws.Range("A18").Value = "Parafusos (Reduzidos)"
$ws.Range("B18").Value = 100

# --- Row 19: TOTAL --- formula stays the same, recalculates automatically
$ws.Range("B19").Formula = "=SUM(B2:B18)"

# --- Page setup / view ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("G21:G24").Select()
